# Remove the "Enterprise Risk management vendors" slide (slide 20), which
# contained a list of external hyperlinks/references to third-party vendor
# sites. Removing it shifts the following slide ("What is what?", sldId 276)
# up into position 20.

$p = $ppt.ActivePresentation

$s = $p.Slides.Item(20)
$s.Delete()
